$d = $word.ActiveDocument

# --- 1. Rename the second "Check" alternate-flow heading to "Fold" ---------
# There are two identical "Check" headings in the document: the first is
# followed by "If the player betted 0 it will count as a if he checked."
# (paragraph 11) and stays untouched; the second is followed by
# "If the player betted -1 it will count as a if he folded." (paragraph 12)
# and is the one that becomes "Fold".
$foldPara = $d.Paragraphs(12)
$foldRange = $foldPara.Range
$null = $foldRange.Find.Execute("Check", $true, $true, $false, $false, $false, `
                                 $true, 1, $false, "Fold", 2)

# --- 2. Move the "_GoBack" bookmark -----------------------------------
# It used to sit right after "Use case: Betting"; now it should sit right
# after the new "Fold" run (before the line break that precedes the
# explanatory sentence).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$foldEnd = $foldRange.End
$bookmarkRange = $d.Range($foldEnd, $foldEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
